$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Relabel the header row (A1:J1) and re-order the rate/balance columns ---
# Columns A-C keep their logical position but get new (Chinese) header text;
# columns D-J are re-ordered together with their relabelling.
$ws.Range("A1").Value = "戶號"
$ws.Range("B1").Value = "額度"
$ws.Range("C1").Value = "撥款序號"
$ws.Range("D1").Value = "下次調息日期"
$ws.Range("E1").Value = "首次調息日期"
$ws.Range("F1").Value = "基本利率代碼"
$ws.Range("G1").Value = "利率加減碼"
$ws.Range("H1").Value = "加碼生效日期"
$ws.Range("I1").Value = "加碼利率"
$ws.Range("J1").Value = "放款餘額"

# --- Widen columns D:I to fit the new (longer) Chinese headers ---
$ws.Columns(4).ColumnWidth = 13.727120535714286
$ws.Columns(5).ColumnWidth = 13.617745535714286
$ws.Columns(6).ColumnWidth = 14.395089285714286
$ws.Columns(7).ColumnWidth = 11.840401785714286
$ws.Columns(8).ColumnWidth = 13.727120535714286
$ws.Columns(9).ColumnWidth = 9.285714285714286

# --- Move the active selection to F4 ---
$ws.Range("F4").Select() | Out-Null
